# horarios_2023_1.xlsx: rename header columns to short names and drop the
# " - PORTUGUÊS E ESPANHOL" suffix from the "LETRAS" course name wherever it
# appears in the "curso" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shorten the header row (row 1): cod_ccr -> cod, cursos_turma -> curso,
#    fases_turma -> fase, expressao_horario -> horario. nome_ccr / ch_ccr
#    (B1 / C1) are left untouched.
$ws.Range("A1").Value = "cod"
$ws.Range("D1").Value = "curso"
$ws.Range("E1").Value = "fase"
$ws.Range("F1").Value = "horario"

# 2) "LETRAS - PORTUGUÊS E ESPANHOL" -> "LETRAS" everywhere it shows up in
#    the curso column (standalone, and as part of "X; LETRAS - PORTUGUÊS E
#    ESPANHOL" combos).
[void]$ws.Columns.Item(4).Replace(" - PORTUGUÊS E ESPANHOL", "")

# 3) Restore the cursor/selection to where the author left it.
[void]$ws.Range("D17").Select()
